$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.194.81'
$ws.Range("E2").Value = '  +6.52%  '
$ws.Range("D3").Value = '3.009.65'
$ws.Range("E3").Value = '  +3.49%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.34'
$ws.Range("E5").Value = '  +2.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.71'
$ws.Range("E6").Value = '  +13.03%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.519'
$ws.Range("E8").Value = '  +3.78%  '
$ws.Range("D9").Value = '3.004.50'
$ws.Range("E9").Value = '  +3.42%  '
$ws.Range("E10").Value = '  -4.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.155'
$ws.Range("E11").Value = '  +4.00%  '
$ws.Range("E12").Value = '  +5.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").Value = '  +6.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.63'
$ws.Range("E14").Value = '  +6.50%  '
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D16").Value = '66.201.90'
$ws.Range("E16").Value = '  +6.67%  '
$ws.Range("D17").Value = '3.513.09'
$ws.Range("E17").Value = '  +3.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.94'
$ws.Range("E18").Value = '  +5.69%  '
$ws.Range("D19").Value = '3.010.34'
$ws.Range("E19").Value = '  +3.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '455.19'
$ws.Range("E20").Value = '  +5.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.94'
$ws.Range("E21").Value = '  +6.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.688'
$ws.Range("E22").Value = '  +5.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.35'
$ws.Range("E23").Value = '  +7.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.37'
$ws.Range("E24").Value = '  +4.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.29'
$ws.Range("E25").Value = '  +14.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.30'
$ws.Range("E26").Value = '  +2.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.54'
$ws.Range("E27").Value = '  +4.83%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.11'
$ws.Range("E29").Value = '  +17.10%  '
$ws.Range("E30").Value = '  +19.07%  '
$ws.Range("E31").Value = '  -6.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.61'
$ws.Range("E32").Value = '  +4.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.27'
$ws.Range("E33").Value = '  +6.25%  '
$ws.Range("E34").Value = '  +5.22%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.996'
$ws.Range("E36").Value = '  +3.96%  '
$ws.Range("E37").Value = '  +16.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.80'
$ws.Range("E38").Value = '  +7.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.02'
$ws.Range("E39").Value = '  +2.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.99'
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.308'
$ws.Range("E41").Value = '  +16.16%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.123'
$ws.Range("E42").Value = '  +8.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.87'
$ws.Range("E43").Value = '  +7.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.45'
$ws.Range("E44").Value = '  +3.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '397.42'
$ws.Range("E45").Value = '  +15.20%  '
$ws.Range("E46").Value = '  +7.07%  '
$ws.Range("D47").Value = '2.800.61'
$ws.Range("E47").Value = '  +2.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.95'
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("E50").Value = '  +11.72%  '
$ws.Range("E51").Value = '  -6.23%  '
